# Change the table style on the three data tables (slides 14, 15 and 16)
# from the old custom style {FD92169A-3642-40FE-9D63-56664CA3861A}
# to {FF8EFB9E-EB65-44FA-AC7C-C2C7635F600C}.

$p = $ppt.ActivePresentation
$oldStyle = "{FD92169A-3642-40FE-9D63-56664CA3861A}"
$newStyle = "{FF8EFB9E-EB65-44FA-AC7C-C2C7635F600C}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyle) {
                $table.ApplyStyle($newStyle)
            }
        }
    }
}
